$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column widths (A:X), ColumnWidth (char units) -> stored width adds 5/6 padding ---
$widths = @(10,12,11,14,14,14,14,104,18,51,9,26,40,27,13,13,31,27,18,33,31,24,76,28)
for ($col = 1; $col -le $widths.Length; $col++) {
    $ws.Columns.Item($col).ColumnWidth = ($widths[$col - 1] - 0.8333333333333334)
}

# --- Header row (A1:W1): horizontal-center + vertical-center + wrap text ---
for ($col = 1; $col -le 23; $col++) {
    $hs = $ws.Cells.Item(1, $col).Style
    $hs.HorizontalAlignment = -4108
    $hs.VerticalAlignment = -4108
    $hs.WrapText = $true
}

# --- X1: same alignment as header, plus yellow fill + updated label ---
$x1s = $ws.Cells.Item(1, 24).Style
$x1s.HorizontalAlignment = -4108
$x1s.VerticalAlignment = -4108
$x1s.WrapText = $true
$x1s.Interior.Color = 65535
$x1s.Interior.PatternColor = 65535
$ws.Range("X1").Value = "Status as of July 11, 2025"

# --- X2: yellow fill only ---
$x2s = $ws.Cells.Item(2, 24).Style
$x2s.Interior.Color = 65535
$x2s.Interior.PatternColor = 65535

# --- Freeze header row ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
